$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.101.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.128.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.28%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "635.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.392"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.767"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.125.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.556"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.30%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.053.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.709.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.135.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000224"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.291.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.156"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.980"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.151"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "504.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.363"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.130"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "163.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0644"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.720"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.74%  "
